$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (Price) values are treated as literal text, not auto-converted numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.282.52'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.921.70'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '0.8079'
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = '244.63'
$ws.Range("E6").Value = '  +0.98%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '0.3277'
$ws.Range("E8").Value = '  +3.55%  '
$ws.Range("D9").Value = '27.09'
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("D10").Value = '0.07284'
$ws.Range("E10").Value = '  +5.38%  '
$ws.Range("D11").Value = '0.7897'
$ws.Range("E11").Value = '  +6.54%  '
$ws.Range("D12").Value = '0.08096'
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("D13").Value = '1.912.98'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").Value = '5.416'
$ws.Range("E14").Value = '  +3.99%  '
$ws.Range("D15").Value = '94.10'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '30.284.97'
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").Value = '14.25'
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").Value = '6.087'
$ws.Range("E18").Value = '  +3.36%  '
$ws.Range("D19").Value = '250.77'
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").Value = '0.000007876'
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("D21").Value = '8.256'
$ws.Range("E21").Value = '  +20.41%  '
$ws.Range("D22").Value = '2.170.73'
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").Value = '0.1640'
$ws.Range("E25").Value = '  +15.03%  '
$ws.Range("D26").Value = '9.494'
$ws.Range("E26").Value = '  +2.72%  '
$ws.Range("D27").Value = '167.97'
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = '19.03'
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = '2.162'
$ws.Range("E29").Value = '  +6.06%  '
$ws.Range("D30").Value = '1.393'
$ws.Range("E30").Value = '  +2.28%  '
$ws.Range("E31").Value = '  +2.38%  '
$ws.Range("D32").Value = '4.409'
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").Value = '0.05713'
$ws.Range("E33").Value = '  +4.16%  '
$ws.Range("D34").Value = '4.152'
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("D35").Value = '1.300'
$ws.Range("E35").Value = '  +2.70%  '
$ws.Range("D36").Value = '0.7511'
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("D37").Value = '1.002'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = '2.730'
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("D39").Value = '0.01961'
$ws.Range("E39").Value = '  +1.75%  '
$ws.Range("D40").Value = '2.828'
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("D41").Value = '0.4547'
$ws.Range("E41").Value = '  +2.66%  '
$ws.Range("D42").Value = '74.33'
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").Value = '6.030'
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").Value = '0.8588'
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").Value = '1.935'
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.037.27'
$ws.Range("E47").Value = '  +5.23%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '103.62'
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("B49").Value = 'SynthetixNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D49").Value = '3.146'
$ws.Range("E49").Value = '  +12.63%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '10.07'
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '7.643'
$ws.Range("E51").Value = '  +1.04%  '

# Restore default style (remove the temporary text-format style) on the Price column
$ws.Range("D2:D51").Style = "Normal"
